# Update "想去人数" (want-to-go count, column F) figures on the
# "展览" (Exhibitions), "演出" (Performances) and "全部类型" (All types)
# sheets to the freshly scraped counts.
$wb = $excel.ActiveWorkbook

# -- 展览 (Exhibitions) --------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 88
$ws1.Range("F4").Value  = 270
$ws1.Range("F6").Value  = 10124
$ws1.Range("F8").Value  = 920
$ws1.Range("F10").Value = 6080
$ws1.Range("F11").Value = 13
$ws1.Range("F12").Value = 418
$ws1.Range("F13").Value = 187
$ws1.Range("F15").Value = 3114
$ws1.Range("F18").Value = 607
$ws1.Range("F20").Value = 23
$ws1.Range("F22").Value = 27
$ws1.Range("F23").Value = 1552

# -- 演出 (Performances) --------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 22

# -- 全部类型 (All types) --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 88
$ws4.Range("F3").Value  = 22
$ws4.Range("F5").Value  = 270
$ws4.Range("F7").Value  = 10124
$ws4.Range("F9").Value  = 920
$ws4.Range("F11").Value = 6080
$ws4.Range("F12").Value = 13
$ws4.Range("F13").Value = 418
$ws4.Range("F14").Value = 187
$ws4.Range("F16").Value = 3114
$ws4.Range("F19").Value = 607
$ws4.Range("F21").Value = 23
$ws4.Range("F23").Value = 27
$ws4.Range("F24").Value = 1552
